$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.78"
$ws.Range("E2").Value = "'1.05%"
$ws.Range("D3").Value = "'36.19"
$ws.Range("E3").Value = "'-1.57%"
$ws.Range("D4").Value = "'5.051"
$ws.Range("E4").Value = "'0.86%"
$ws.Range("D5").Value = "'0.07899"
$ws.Range("E5").Value = "'2.20%"
$ws.Range("D6").Value = "'2.290"
$ws.Range("E6").Value = "'9.57%"
$ws.Range("D7").Value = "'7.996"
$ws.Range("E7").Value = "'0.02%"
$ws.Range("D8").Value = "'4.151"
$ws.Range("E8").Value = "'2.47%"
$ws.Range("D9").Value = "'0.9285"
$ws.Range("E9").Value = "'0.98%"
$ws.Range("D10").Value = "'0.09852"
$ws.Range("E10").Value = "'0.85%"
$ws.Range("D11").Value = "'0.1869"
$ws.Range("E11").Value = "'0.53%"
$ws.Range("D12").Value = "'0.08962"
$ws.Range("E12").Value = "'3.74%"
$ws.Range("D13").Value = "'0.03754"
$ws.Range("E13").Value = "'4.39%"
$ws.Range("D14").Value = "'0.09917"
$ws.Range("E14").Value = "'-0.50%"
$ws.Range("D15").Value = "'0.001442"
$ws.Range("E15").Value = "'-2.03%"
$ws.Range("D16").Value = "'0.005718"
$ws.Range("E16").Value = "'-0.40%"
$ws.Range("D17").Value = "'3.466"
$ws.Range("E17").Value = "'-0.20%"
$ws.Range("D18").Value = "'2.632"
$ws.Range("E18").Value = "'3.07%"
$ws.Range("D19").Value = "'0.3367"
$ws.Range("E19").Value = "'-1.92%"
$ws.Range("D20").Value = "'0.1319"
$ws.Range("E20").Value = "'-0.99%"
$ws.Range("D21").Value = "'5.075"
$ws.Range("E21").Value = "'2.27%"
$ws.Range("E22").Value = "'1.60%"
$ws.Range("D23").Value = "'0.04579"
$ws.Range("E23").Value = "'-1.00%"
$ws.Range("D24").Value = "'0.001234"
$ws.Range("E24").Value = "'-0.34%"
$ws.Range("D25").Value = "'0.004777"
$ws.Range("E25").Value = "'-6.51%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'-7.62%"
$ws.Range("D39").Value = "'0.01924"
$ws.Range("E39").Value = "'8.48%"
$ws.Range("D40").Value = "'0.04916"
$ws.Range("E40").Value = "'5.54%"
$ws.Range("D41").Value = "'0.007869"
$ws.Range("E41").Value = "'2.25%"
$ws.Range("D42").Value = "'0.1391"
$ws.Range("E42").Value = "'0.02%"
$ws.Range("E43").Value = "'-2.08%"
$ws.Range("D44").Value = "'0.002183"
$ws.Range("E44").Value = "'-3.18%"
$ws.Range("E45").Value = "'15.44%"
$ws.Range("D46").Value = "'0.00006147"
$ws.Range("E46").Value = "'-2.06%"
$ws.Range("E47").Value = "'-0.56%"
$ws.Range("E48").Value = "'51.36%"
$ws.Range("D49").Value = "'0.001802"
$ws.Range("E49").Value = "'-10.39%"
$ws.Range("E50").Value = "'-0.56%"
$ws.Range("E51").Value = "'-0.56%"
